$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Wednesday, January 1, 2020 00:00:00"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
